$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly price data refresh: row contents (date + volume/price fields) were
# reshuffled across rows 2-30 (dates/prices re-assigned to different weekly
# observations). Row 19 is unchanged.

# Row 2 <- data from original row 30
$ws.Range("D2").Value = 44257
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 1400
$ws.Range("M2").Value = 1450
$ws.Range("P2").Value = 725

# Row 3 <- data from original row 23
$ws.Range("D3").Value = 44363
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 2800
$ws.Range("M3").Value = 2650
$ws.Range("P3").Value = 1325

# Row 4 <- data from original row 3
$ws.Range("D4").Value = 44540

# Row 5 <- data from original row 9
$ws.Range("D5").Value = 44447
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 900
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = 950
$ws.Range("P5").Value = 475

# Row 6 <- data from original row 14
$ws.Range("D6").Value = 44243
$ws.Range("K6").Value = 1200
$ws.Range("L6").Value = 1300
$ws.Range("M6").Value = 1250
$ws.Range("P6").Value = 625

# Row 7 <- data from original row 17
$ws.Range("D7").Value = 44385
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 2400
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = 2450
$ws.Range("P7").Value = 1225

# Row 8 <- data from original row 27
$ws.Range("D8").Value = 44525
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 1400
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = 1450
$ws.Range("P8").Value = 725

# Row 9 <- data from original row 4
$ws.Range("D9").Value = 44468

# Row 10 <- data from original row 25
$ws.Range("D10").Value = 44302
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 900
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = 950
$ws.Range("P10").Value = 475

# Row 11 <- data from original row 5
$ws.Range("D11").Value = 44291
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 1800
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 1900
$ws.Range("P11").Value = 950

# Row 12 <- data from original row 24
$ws.Range("D12").Value = 44435
$ws.Range("K12").Value = 900
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = 950
$ws.Range("P12").Value = 475

# Row 13 <- data from original row 10
$ws.Range("D13").Value = 44390
$ws.Range("J13").Value = 250
$ws.Range("K13").Value = 2400
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = 2450
$ws.Range("P13").Value = 1225

# Row 14 <- data from original row 11
$ws.Range("D14").Value = 44601
$ws.Range("J14").Value = 270
$ws.Range("K14").Value = 2200
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2350
$ws.Range("P14").Value = 1175

# Row 15 <- data from original row 16
$ws.Range("D15").Value = 44544
$ws.Range("K15").Value = 900
$ws.Range("L15").Value = 1000
$ws.Range("M15").Value = 950
$ws.Range("P15").Value = 475

# Row 16 <- data from original row 13
$ws.Range("D16").Value = 44161
$ws.Range("J16").Value = 270

# Row 17 <- data from original row 28
$ws.Range("D17").Value = 44365
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 1800
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = 1900
$ws.Range("P17").Value = 950

# Row 18 <- data from original row 7
$ws.Range("D18").Value = 44403
$ws.Range("J18").Value = 250
$ws.Range("K18").Value = 1800
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = 1900
$ws.Range("P18").Value = 950

# Row 20 <- data from original row 8
$ws.Range("D20").Value = 44392

# Row 21 <- data from original row 2
$ws.Range("D21").Value = 44172
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 1300
$ws.Range("L21").Value = 1500
$ws.Range("M21").Value = 1400
$ws.Range("P21").Value = 700

# Row 22 <- data from original row 20
$ws.Range("D22").Value = 44202

# Row 23 <- data from original row 15
$ws.Range("D23").Value = 44726

# Row 24 <- data from original row 6
$ws.Range("D24").Value = 44427
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = 1300
$ws.Range("L24").Value = 1500
$ws.Range("M24").Value = 1400
$ws.Range("P24").Value = 700

# Row 25 <- data from original row 22
$ws.Range("D25").Value = 44253
$ws.Range("J25").Value = 250
$ws.Range("K25").Value = 1800
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = 1900
$ws.Range("P25").Value = 950

# Row 26 <- data from original row 12
$ws.Range("D26").Value = 44266
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 1700
$ws.Range("L26").Value = 1800
$ws.Range("M26").Value = 1750
$ws.Range("N26").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P26").Value = 875
$ws.Range("Q26").Value = 2

# Row 27 <- data from original row 29
$ws.Range("D27").Value = 44229
$ws.Range("J27").Value = 250
$ws.Range("K27").Value = 1800
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = 1900
$ws.Range("P27").Value = 950

# Row 28 <- data from original row 26
$ws.Range("D28").Value = 44181
$ws.Range("K28").Value = 1000
$ws.Range("L28").Value = 1200
$ws.Range("M28").Value = 1100
$ws.Range("N28").Value = '$/atado'
$ws.Range("P28").Value = 1100
$ws.Range("Q28").Value = 1

# Row 29 <- data from original row 21
$ws.Range("D29").Value = 44438
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 950
$ws.Range("L29").Value = 1000
$ws.Range("M29").Value = 975
$ws.Range("P29").Value = 488

# Row 30 <- data from original row 18
$ws.Range("D30").Value = 44572
$ws.Range("J30").Value = 300
